$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 08:09"

# --- 2) Swap adjacent country-name pairs (shared-string reorder in the
#        source diff manifests as the displayed country name swapping
#        between two consecutive rows; the numeric stats for those rows
#        are unaffected except where noted separately below) ---

# Fiyi (row 202) <-> Dominica (row 203)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Groenlandia (row 208) <-> Islas Malvinas (row 209)
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

# Seychelles (row 211) <-> Montserrat (row 212)
$ws.Range("A211").Value = "Montserrat"
$ws.Range("A212").Value = "Seychelles"

# --- 3) Numeric data updates ---

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 7228
$ws.Range("C75").Value = 51
$ws.Range("E75").Value = 2331

# Row 97 - Tailandia
$ws.Range("B97").Value = 3162
$ws.Range("C97").Value = 4
$ws.Range("D97").Value = 3040
$ws.Range("E97").Value = 64

# Row 189 - Macao
$ws.Range("B189").Value = 46
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 1

# Row 211 (now showing "Montserrat")
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 212 (now showing "Seychelles")
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
